$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet and update workbook title date
$ws.Name = "Through 2022-09-12"

# Update the "September (through 09-11)" label cell to reflect new date
$ws.Range("A10").Value = "September (through 09-12)"

# Update September row (row 10) values
$ws.Range("B10").Value = 11
$ws.Range("C10").Value = 20
$ws.Range("D10").Value = 30
$ws.Range("E10").Value = 24
$ws.Range("F10").Value = 28
$ws.Range("G10").Value = 39
$ws.Range("H10").Value = 60
$ws.Range("I10").Value = 61

# Update Total row (row 11) values
$ws.Range("B11").Value = 205
$ws.Range("C11").Value = 401
$ws.Range("D11").Value = 581
$ws.Range("E11").Value = 514
$ws.Range("F11").Value = 377
$ws.Range("G11").Value = 823
$ws.Range("H11").Value = 1130
$ws.Range("I11").Value = 1198
